$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.337.30"
$ws.Range("E2").Value = "  -0.16%  "
Set-TextValue $ws.Range("D3") "1.841.11"
$ws.Range("E3").Value = "  -0.60%  "
Set-TextValue $ws.Range("D4") "0.9987"
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue $ws.Range("D5") "240.20"
$ws.Range("E5").Value = "  -0.28%  "
Set-TextValue $ws.Range("D6") "0.6280"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E7").Value = "  +0.08%  "
Set-TextValue $ws.Range("D8") "0.07421"
$ws.Range("E8").Value = "  -2.19%  "
Set-TextValue $ws.Range("D9") "0.2896"
$ws.Range("E9").Value = "  -1.02%  "
Set-TextValue $ws.Range("D10") "24.82"
$ws.Range("E10").Value = "  +1.22%  "
Set-TextValue $ws.Range("D11") "0.07733"
$ws.Range("E11").Value = "  -0.13%  "
Set-TextValue $ws.Range("D12") "1.842.09"
$ws.Range("E12").Value = "  -0.56%  "
Set-TextValue $ws.Range("D13") "4.978"
$ws.Range("E13").Value = "  -0.96%  "
Set-TextValue $ws.Range("D14") "0.6772"
$ws.Range("E14").Value = "  -0.90%  "
Set-TextValue $ws.Range("D15") "0.00001021"
$ws.Range("E15").Value = "  -2.49%  "
Set-TextValue $ws.Range("D16") "81.95"
$ws.Range("E16").Value = "  -1.68%  "
Set-TextValue $ws.Range("D17") "6.245"
$ws.Range("E17").Value = "  +1.63%  "
Set-TextValue $ws.Range("D18") "29.384.07"
$ws.Range("E18").Value = "  -0.03%  "
Set-TextValue $ws.Range("D19") "229.01"
$ws.Range("E19").Value = "  -0.53%  "
Set-TextValue $ws.Range("D20") "12.30"
$ws.Range("E20").Value = "  -0.70%  "
Set-TextValue $ws.Range("D21") "0.9996"
$ws.Range("E21").Value = "  +0.06%  "
Set-TextValue $ws.Range("D22") "7.414"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("E23").Value = "  +0.11%  "
Set-TextValue $ws.Range("D24") "158.96"
$ws.Range("E24").Value = "  -0.04%  "
Set-TextValue $ws.Range("D25") "8.471"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -3.39%  "
Set-TextValue $ws.Range("D27") "17.41"
$ws.Range("E27").Value = "  -1.40%  "
Set-TextValue $ws.Range("D28") "0.06510"
$ws.Range("E28").Value = "  +14.37%  "
Set-TextValue $ws.Range("D29") "1.439"
$ws.Range("E29").Value = "  +1.53%  "
Set-TextValue $ws.Range("D30") "1.485"
$ws.Range("E30").Value = "  +0.52%  "
Set-TextValue $ws.Range("D31") "4.065"
$ws.Range("E31").Value = "  -1.71%  "
Set-TextValue $ws.Range("D32") "4.070"
$ws.Range("E32").Value = "  +0.29%  "
Set-TextValue $ws.Range("D33") "1.834"
$ws.Range("E33").Value = "  +0.32%  "
Set-TextValue $ws.Range("D34") "1.138"
$ws.Range("E34").Value = "  -1.64%  "
Set-TextValue $ws.Range("D35") "0.6910"
$ws.Range("E35").Value = "  -1.44%  "
Set-TextValue $ws.Range("D36") "2.560"
$ws.Range("E36").Value = "  -0.75%  "
Set-TextValue $ws.Range("D37") "0.01855"
$ws.Range("E37").Value = "  +1.54%  "
Set-TextValue $ws.Range("D38") "2.819"
$ws.Range("E38").Value = "  +3.35%  "
Set-TextValue $ws.Range("D39") "1.243.97"
$ws.Range("E39").Value = "  -0.22%  "
Set-TextValue $ws.Range("D40") "6.730"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("B42").Value = "RocketPoolETH"
$ws.Range("C42").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D42") "2.027.64"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D43") "0.9995"
$ws.Range("E43").Value = "  +0.04%  "
Set-TextValue $ws.Range("D44") "100.79"
$ws.Range("E44").Value = "  -1.31%  "
Set-TextValue $ws.Range("D45") "65.64"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D46") "0.00000000119"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D47") "7.050"
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D48") "1.714"
$ws.Range("E48").Value = "  +2.18%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D49") "0.1151"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "8.984"
$ws.Range("E50").Value = "  -0.77%  "
Set-TextValue $ws.Range("D51") "0.3890"
$ws.Range("E51").Value = "  -2.00%  "
